# Update w/ harris poll (12/31) and ow rolling poll (1/6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 150 (opinionway rolling poll gets corrected/updated values) ---
$ws.Range("C150").Value = 19
$ws.Range("E150").Value = 3
$ws.Range("I150").Value = 761

# --- Add new row 152: opinionway online rolling poll, week 19, 2022-01-05 ---
$ws.Range("A152").Value = 58
$ws.Range("B152").Value = 2022
$ws.Range("C152").Value = 19
$ws.Range("D152").Value = 1
$ws.Range("E152").Value = 5
$ws.Range("F152").Value = "opinionway"
$ws.Range("G152").Value = "online"
$ws.Range("H152").Value = "partially"
$ws.Range("I152").Value = 1534
$ws.Range("J152").Value = 1
$ws.Range("K152").Value = 1
$ws.Range("L152").Value = 9
$ws.Range("M152").Value = 3
$ws.Range("N152").Value = 1
$ws.Range("O152").Value = 7
$ws.Range("P152").Value = 4
$ws.Range("Q152").Value = 25
$ws.Range("R152").Value = 17
$ws.Range("U152").Value = 1
$ws.Range("V152").Value = 2
$ws.Range("W152").Value = 17
$ws.Range("X152").Value = 12

# --- Add new row 153: harris online poll, week 18, 2021-12-30 ---
$ws.Range("A153").Value = 59
$ws.Range("B153").Value = 2021
$ws.Range("C153").Value = 18
$ws.Range("D153").Value = 12
$ws.Range("E153").Value = 30
$ws.Range("F153").Value = "harris"
$ws.Range("G153").Value = "online"
$ws.Range("H153").Value = "included"
$ws.Range("I153").Value = 2183
$ws.Range("J153").Value = 1
$ws.Range("K153").Value = "T_0.5"
$ws.Range("L153").Value = 10
$ws.Range("M153").Value = 3
$ws.Range("N153").Value = 1
$ws.Range("O153").Value = 7
$ws.Range("P153").Value = 4
$ws.Range("Q153").Value = 24
$ws.Range("R153").Value = 16
$ws.Range("U153").Value = "T_0.5"
$ws.Range("V153").Value = 2
$ws.Range("W153").Value = 16
$ws.Range("X153").Value = 16
$ws.Range("Y153").Value = "T_0.5"
$ws.Range("AA153").Value = "T_0.5"

# --- Update the view state to match the new scroll/selection position ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 132
$win.ScrollColumn = 7
$ws.Range("Y154").Select() | Out-Null
